$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (column D) and 1h volume change (column E) values,
# refreshed by the "Updated cryptos list" GitHub Actions job.
$updates = @{
    'D2' = '69.831.17'
    'E2' = '  +4.25%  '
    'D3' = '3.630.84'
    'E4' = '  -0.11%  '
    'D5' = '632.12'
    'E5' = '  +3.59%  '
    'D6' = '160.06'
    'E6' = '  +5.33%  '
    'D7' = '3.629.35'
    'E7' = '  +3.21%  '
    'E8' = '  -0.08%  '
    'D9' = '0.497'
    'E9' = '  +2.84%  '
    'E10' = '  +6.42%  '
    'D11' = '7.34'
    'E11' = '  +6.70%  '
    'D12' = '0.442'
    'E12' = '  +3.55%  '
    'E13' = '  +4.76%  '
    'E14' = '  +6.23%  '
    'D15' = '4.246.37'
    'E15' = '  +3.22%  '
    'D16' = '3.634.87'
    'E16' = '  +3.41%  '
    'D17' = '69.858.41'
    'E17' = '  +4.28%  '
    'E18' = '  +0.29%  '
    'D19' = '6.69'
    'E19' = '  +6.45%  '
    'D20' = '16.09'
    'E20' = '  +5.03%  '
    'D21' = '10.13'
    'E21' = '  +11.10%  '
    'D22' = '465.04'
    'E22' = '  +4.74%  '
    'E23' = '  +2.80%  '
    'D24' = '78.99'
    'E24' = '  +1.72%  '
    'E25' = '  +12.80%  '
    'D26' = '10.83'
    'E26' = '  +6.21%  '
    'D27' = '3.776.70'
    'E27' = '  +3.21%  '
    'E28' = '  +0.03%  '
    'D29' = '9.38'
    'E29' = '  +15.22%  '
    'E30' = '  +4.88%  '
    'D31' = '1.73'
    'E31' = '  +5.56%  '
    'D33' = '6.64'
    'E33' = '  +7.95%  '
    'D34' = '0.999'
    'E34' = '  +0.03%  '
    'E35' = '  +6.43%  '
    'D36' = '26.66'
    'E36' = '  +3.85%  '
    'D37' = '3.629.55'
    'E37' = '  +3.46%  '
    'E38' = '  +5.48%  '
    'E39' = '  +14.04%  '
    'E40' = '  +0.01%  '
    'E41' = '  +8.69%  '
    'D42' = '178.89'
    'E42' = '  +3.32%  '
    'D43' = '0.999'
    'E43' = '  -0.10%  '
    'E44' = '  +1.98%  '
    'D45' = '31.61'
    'E45' = '  +16.79%  '
    'D46' = '0.916'
    'E46' = '  +3.22%  '
    'D47' = '1.37'
    'E47' = '  +12.71%  '
    'E48' = '  +10.31%  '
    'D49' = '46.47'
    'E49' = '  +2.81%  '
    'E50' = '  +3.90%  '
    'E51' = '  +9.47%  '
}

# Cells whose new value reads as a plain number (e.g. "632.12"). These must
# be kept as text, matching the original inline-string cell type, so force
# text formatting before writing the value to stop Excel auto-converting
# them to a numeric type.
$forceTextCells = @(
    'D5','D6','D9','D11','D12','D19','D20','D21','D22','D24','D26','D29',
    'D31','D33','D34','D36','D42','D43','D45','D46','D47','D49'
)

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    if ($forceTextCells -contains $cellRef) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$cellRef]
}
